# Add new combat tactics, balancing of existing, new combat phase
#
# This adds a new "Sheet2" worksheet (right after "Sheet1") that holds the
# tactics for the new "recon" combat phase: a small Attack/Defense
# cross-reference table plus the counter lists for each side.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet directly after Sheet1 so sheet order/relationship
# ids come out as Sheet1, Sheet2 (matches workbook.xml <sheets> order).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Column E: list of new "recon" attack-side tactics -------------------
# (typed first, top to bottom)
$ws2.Range("E12").Value = "Limited recon"
$ws2.Range("E13").Value = "Aggressive recon"
$ws2.Range("E14").Value = "Recon Attack"
$ws2.Range("E15").Value = "Supression advance"
$ws2.Range("E16").Value = "Mobile recon"
$ws2.Range("E17").Value = "breakthrough recon"

# --- Header row -----------------------------------------------------------
$ws2.Range("E11").Value = "Attack"
$ws2.Range("F11").Value = "Defense"

# --- Column F: matching defense-side counters ------------------------------
$ws2.Range("F14").Value = "Recon defense"
$ws2.Range("F12").Value = "Hunter packs"
$ws2.Range("F13").Value = "Ambush formation"
$ws2.Range("F15").Value = "Counter battery"
$ws2.Range("F16").Value = "Node defense"

# --- Trailing note row ------------------------------------------------------
$ws2.Range("E18").Value = "Pincer manouver"

# Column widths for the new sheet.
$ws2.Range("E1").ColumnWidth = 24.833333333333332
$ws2.Range("F1").ColumnWidth = 24.666666666666668

# Leave the selection on Sheet2 where the last edits were made, then hop
# back to Sheet1 (which stays the active/visible tab) and restore its
# selection.
[void]$ws2.Range("E17").Select()

[void]$ws1.Activate()
[void]$ws1.Range("C23").Select()
